# This edit re-orders the data rows 2-7 of the "Artfynd" sheet (a cyclic
# permutation of full rows, each row keeping all of its own cell content
# together, including the "På död gran" comment that travels with the row
# that previously lived at row 7).
#
# Mapping (destination row <- source row), using the row numbers as they
# exist in the original/"before" worksheet:
#   row 2 <- row 5
#   row 3 <- row 7
#   row 4 <- row 2
#   row 5 <- row 3
#   row 6 <- row 4
#   row 7 <- row 6
#
# Because several rows both provide data to, and receive data from, other
# rows in this block, the whole 2-7 block is first staged into a scratch
# area of the sheet, and then copied back into its final, re-ordered
# position. Using Range.Copy (rather than reading/writing .Value2 arrays)
# preserves each cell's original data type (numbers stay numbers, text that
# looks like a date such as "2019-09-17" stays text, booleans stay
# booleans) and keeps the existing cell styling untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcFirst = 2
$srcLast = 7
$lastCol = "AY"

# Scratch area used to stage the original rows 2-7 while they are
# rearranged (chosen far below the real data so it cannot collide with it).
$stageFirst = 200

# destination row (sheet row number) -> source row (sheet row number)
$rowMap = @{
    2 = 5
    3 = 7
    4 = 2
    5 = 3
    6 = 4
    7 = 6
}

# 1. Copy rows 2-7 into the scratch area, in original order, rows
#    stageFirst .. stageFirst+5 correspond to original rows 2 .. 7.
$stageLast = $stageFirst + ($srcLast - $srcFirst)
$ws.Range("A" + $stageFirst + ":" + $lastCol + $stageLast).Clear()
$ws.Range("A" + $srcFirst + ":" + $lastCol + $srcLast).Copy($ws.Range("A" + $stageFirst + ":" + $lastCol + $stageLast))

# 2. Clear the original block so that blank source cells (e.g. the removed
#    "På död gran" comment) really end up blank in the destination rather
#    than leaving stale content behind.
$ws.Range("A" + $srcFirst + ":" + $lastCol + $srcLast).Clear()

# 3. Copy each staged row back to its final destination row.
foreach ($destRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$destRow]
    $stageRow = $stageFirst + ($sourceRow - $srcFirst)
    $srcRange = $ws.Range("A" + $stageRow + ":" + $lastCol + $stageRow)
    $dstRange = $ws.Range("A" + $destRow + ":" + $lastCol + $destRow)
    $srcRange.Copy($dstRange)
}

# 4. Clean up the scratch area.
$ws.Range("A" + $stageFirst + ":" + $lastCol + $stageLast).Clear()
